# CPP/CPPbI: switch to new industry categories (#89)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "CPPbI": replace the old 8 broad industry categories with
# the new 25 detailed industry categories
# ---------------------------------------------------------------
$cppbi = $wb.Worksheets.Item("CPPbI")

# Clear old data rows (2-9) first
$cppbi.Range("A2:C9").ClearContents()

$industries = @(
    @("agriculture and forestry 01T03", 0, 0),
    @("coal mining 05", 0, 0),
    @("oil and gas extraction 06", 0, 0),
    @("other mining and quarrying 07T08", 0, 0),
    @("food beverage and tobacco 10T12", 1, 1),
    @("textiles apparel and leather 13T15", 1, 1),
    @("wood products 16", 1, 1),
    @("pulp paper and printing 17T18", 1, 1),
    @("refined petroleum and coke 19", 1, 1),
    @("chemicals 20", 1, 1),
    @("rubber and plastic products 22", 1, 1),
    @("glass and glass products 231", 1, 1),
    @("cement and other nonmetallic minerals 239", 1, 1),
    @("iron and steel 241", 1, 1),
    @("other metals 242", 1, 1),
    @("metal products except machinery and vehicles 25", 1, 1),
    @("computers and electronics 26", 1, 1),
    @("appliances and electrical equipment 27", 1, 1),
    @("other machinery 28", 1, 1),
    @("road vehicles 29", 1, 1),
    @("nonroad vehicles 30", 1, 1),
    @("other manufacturing 31T33", 1, 1),
    @("energy pipelines and gas processing 352T353", 1, 1),
    @("water and waste 36T39", 0, 0),
    @("construction 41T43", 0, 0)
)

$row = 2
foreach ($item in $industries) {
    $cppbi.Cells.Item($row, 1).Value = $item[0]
    $cppbi.Cells.Item($row, 2).Value = $item[1]
    $cppbi.Cells.Item($row, 3).Value = $item[2]
    $row = $row + 1
}

# Column A is slightly wider to fit the longer industry names
$cppbi.Columns.Item(1).ColumnWidth = 45.42578125

# ---------------------------------------------------------------
# Sheet "About": update explanatory notes
# ---------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

# Row 16 note text is updated to mention construction & water/waste
$about.Range("A16").Value = "For industries, we assign 100% to all industries except mining, agriculture, construction, and water and waste, as the activities"

# Row 20 gains a new explanatory note (previously blank)
$about.Range("A20").Value = '(We assume CO2 from "water and waste" is from waste collection trucks, not water treatment plants, which use almost entirely electricity.)'

# Insert a new blank row after row 20 (row 21), pushing the remaining notes down
$about.Rows.Item(21).Insert()
